# Update the worksheet date and each division problem to the new values.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-11 Wednesday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-09-12 Thursday", 2) | Out-Null
$d.Content.Find.Execute("796÷2=", $true, $true, $false, $false, $false, $true, 1, $false, "254÷3=", 2) | Out-Null
$d.Content.Find.Execute("150÷5=", $true, $true, $false, $false, $false, $true, 1, $false, "778÷7=", 2) | Out-Null
$d.Content.Find.Execute("493÷5=", $true, $true, $false, $false, $false, $true, 1, $false, "809÷3=", 2) | Out-Null
$d.Content.Find.Execute("699÷4=", $true, $true, $false, $false, $false, $true, 1, $false, "566÷2=", 2) | Out-Null
$d.Content.Find.Execute("737÷5=", $true, $true, $false, $false, $false, $true, 1, $false, "828÷5=", 2) | Out-Null
$d.Content.Find.Execute("245÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "564÷9=", 2) | Out-Null
$d.Content.Find.Execute("531÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "408÷2=", 2) | Out-Null
$d.Content.Find.Execute("290÷5=", $true, $true, $false, $false, $false, $true, 1, $false, "138÷6=", 2) | Out-Null
$d.Content.Find.Execute("172÷3=", $true, $true, $false, $false, $false, $true, 1, $false, "694÷9=", 2) | Out-Null
$d.Content.Find.Execute("983÷8=", $true, $true, $false, $false, $false, $true, 1, $false, "332÷8=", 2) | Out-Null
$d.Content.Find.Execute("418÷8=", $true, $true, $false, $false, $false, $true, 1, $false, "561÷2=", 2) | Out-Null
$d.Content.Find.Execute("695÷5=", $true, $true, $false, $false, $false, $true, 1, $false, "662÷9=", 2) | Out-Null
$d.Content.Find.Execute("672÷6=", $true, $true, $false, $false, $false, $true, 1, $false, "114÷7=", 2) | Out-Null
$d.Content.Find.Execute("573÷6=", $true, $true, $false, $false, $false, $true, 1, $false, "733÷6=", 2) | Out-Null
$d.Content.Find.Execute("506÷8=", $true, $true, $false, $false, $false, $true, 1, $false, "307÷4=", 2) | Out-Null
$d.Content.Find.Execute("451÷9=", $true, $true, $false, $false, $false, $true, 1, $false, "882÷2=", 2) | Out-Null
$d.Content.Find.Execute("569÷5=", $true, $true, $false, $false, $false, $true, 1, $false, "401÷8=", 2) | Out-Null
$d.Content.Find.Execute("297÷8=", $true, $true, $false, $false, $false, $true, 1, $false, "215÷6=", 2) | Out-Null
$d.Content.Find.Execute("402÷2=", $true, $true, $false, $false, $false, $true, 1, $false, "117÷9=", 2) | Out-Null
$d.Content.Find.Execute("728÷5=", $true, $true, $false, $false, $false, $true, 1, $false, "147÷3=", 2) | Out-Null
$d.Content.Find.Execute("301÷2=", $true, $true, $false, $false, $false, $true, 1, $false, "281÷4=", 2) | Out-Null
$d.Content.Find.Execute("307÷9=", $true, $true, $false, $false, $false, $true, 1, $false, "955÷3=", 2) | Out-Null
$d.Content.Find.Execute("265÷8=", $true, $true, $false, $false, $false, $true, 1, $false, "902÷2=", 2) | Out-Null
$d.Content.Find.Execute("426÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "587÷7=", 2) | Out-Null
$d.Content.Find.Execute("104÷6=", $true, $true, $false, $false, $false, $true, 1, $false, "726÷2=", 2) | Out-Null
